$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Insert two new rows before the current row 7 ("Quận/Huyện:" row), pushing
# everything below (rows 7-15) down to rows 9-17. The new rows inherit the
# formatting of the row above (row 6): style 2 on column D, style 8 on E:G.
$ws.Rows("7:8").Insert()

# --- Filter / search block (rows 5-13) -------------------------------------
# Row 5: add an empty B5 cell (format copied from B4) and relabel D5
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D5").Value = "Từ khóa:"

# Row 6: add an empty B6 cell (format copied from B4) and relabel D6
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D6").Value = "Người nộp đơn:"

$excel.CutCopyMode = $false

# Row 7 (newly inserted): "Lĩnh vực:"
$ws.Range("D7").Value = "Lĩnh vực:"

# Row 8 (newly inserted): "Tỉnh/TP:"
$ws.Range("D8").Value = "Tỉnh/TP:"

# Row 9 (was row 7): "Quận/Huyện:"
$ws.Range("D9").Value = "Quận/Huyện:"

# Row 10 (was row 8): "Xã/Phường:"
$ws.Range("D10").Value = "Xã/Phường:"

# Row 11 (was row 9): "Thời gian tiếp nhận:" (fixes the stray leading "t")
$ws.Range("D11").Value = "Thời gian tiếp nhận:"

# Row 12 (was row 10): "Kết quả:"
$ws.Range("D12").Value = "Kết quả:"

# Row 13 (was row 11): "Công khai:"
$ws.Range("D13").Value = "Công khai:"

# Header row (was row 14, now row 16) and data row (was row 15, now row 17)
# keep their original content automatically thanks to the row insert shift.

# --- Refresh the frozen pane / selection to match the new row layout -------
$win.FreezePanes = $false
$ws.Range("A17").Select()
$win.FreezePanes = $true
$ws.Range("D13").Select()

$wb.Save()
